$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 25,14
$data[0,0] = 1.521407
$data[0,1] = 4.564221
$data[0,2] = 0.07220977817024694
$data[0,3] = 0.07583348915219905
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 3.721411333333334
$data[0,7] = 11.164234
$data[0,8] = 0.1546085275451463
$data[0,9] = 0.1613881693704714
$data[0,10] = 5.661781252412667
$data[0,11] = 50.956031271714
$data[0,12] = 0.01116424747726353
$data[0,13] = 0.0122386279912489
$data[1,0] = 1.521407
$data[1,1] = 4.564221
$data[1,2] = 0.07220977817024694
$data[1,3] = 0.07583348915219905
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 5.194472666666667
$data[1,7] = 15.583418
$data[1,8] = 0.2158078477305768
$data[1,9] = 0.2252711026618443
$data[1,10] = 7.902907076375334
$data[1,11] = 71.126163687378
$data[1,12] = 0.01558343681202338
$data[1,13] = 0.01708309372001088
$data[2,0] = 1.521407
$data[2,1] = 4.564221
$data[2,2] = 0.07220977817024694
$data[2,3] = 0.07583348915219905
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 5.602208666666667
$data[2,7] = 16.806626
$data[2,8] = 0.2327475130727259
$data[2,9] = 0.2429535786722285
$data[2,10] = 8.523239480927334
$data[2,11] = 76.709155328346
$data[2,12] = 0.01680664628865819
$data[2,13] = 0.01842401757272838
$data[3,0] = 1.521407
$data[3,1] = 4.564221
$data[3,2] = 0.07220977817024694
$data[3,3] = 0.07583348915219905
$data[3,4] = 2
$data[3,5] = 1
$data[3,6] = 3.033406
$data[3,7] = 6.066812
$data[3,8] = 0.1260248849423826
$data[3,9] = 0.08770074889103974
$data[3,10] = 4.615045122242
$data[3,11] = 27.690270733452
$data[3,12] = 0.00910022898562034
$data[3,13] = 0.006650653789668394
$data[4,0] = 1.521407
$data[4,1] = 4.564221
$data[4,2] = 0.07220977817024694
$data[4,3] = 0.07583348915219905
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 6.518398333333334
$data[4,7] = 19.555195
$data[4,8] = 0.2708112267091685
$data[4,9] = 0.2826864004044161
$data[4,10] = 9.917136853121667
$data[4,11] = 89.254231678095
$data[4,12] = 0.01955521860668151
$data[4,13] = 0.02143709607854248
$data[5,0] = 4.966822333333333
$data[5,1] = 14.900467
$data[5,2] = 0.2357377998793408
$data[5,3] = 0.2475678549761722
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 3.721411333333334
$data[5,7] = 11.164234
$data[5,8] = 0.1546085275451463
$data[5,9] = 0.1613881693704714
$data[5,10] = 18.48358892191978
$data[5,11] = 166.352300297278
$data[5,12] = 0.03644707412607724
$data[5,13] = 0.03995452290957878
$data[6,0] = 4.966822333333333
$data[6,1] = 14.900467
$data[6,2] = 0.2357377998793408
$data[6,3] = 0.2475678549761722
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 5.194472666666667
$data[6,7] = 15.583418
$data[6,8] = 0.2158078477305768
$data[6,9] = 0.2252711026618443
$data[6,10] = 25.80002285068955
$data[6,11] = 232.200205656206
$data[6,12] = 0.05087406722070196
$data[6,13] = 0.05576988367410987
$data[7,0] = 4.966822333333333
$data[7,1] = 14.900467
$data[7,2] = 0.2357377998793408
$data[7,3] = 0.2475678549761722
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 5.602208666666667
$data[7,7] = 16.806626
$data[7,8] = 0.2327475130727259
$data[7,9] = 0.2429535786722285
$data[7,10] = 27.82517512159356
$data[7,11] = 250.426576094342
$data[7,12] = 0.05486738665915253
$data[7,13] = 0.06014749633066831
$data[8,0] = 4.966822333333333
$data[8,1] = 14.900467
$data[8,2] = 0.2357377998793408
$data[8,3] = 0.2475678549761722
$data[8,4] = 2
$data[8,5] = 1
$data[8,6] = 3.033406
$data[8,7] = 6.066812
$data[8,8] = 0.1260248849423826
$data[8,9] = 0.08770074889103974
$data[8,10] = 15.06638866686733
$data[8,11] = 90.39833200120398
$data[8,12] = 0.02970882910636434
$data[8,13] = 0.02171188628275862
$data[9,0] = 4.966822333333333
$data[9,1] = 14.900467
$data[9,2] = 0.2357377998793408
$data[9,3] = 0.2475678549761722
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 6.518398333333334
$data[9,7] = 19.555195
$data[9,8] = 0.2708112267091685
$data[9,9] = 0.2826864004044161
$data[9,10] = 32.37572641956278
$data[9,11] = 291.381537776065
$data[9,12] = 0.06384044276704476
$data[9,13] = 0.06998406577905664
$data[10,0] = 4.707986666666667
$data[10,1] = 14.12396
$data[10,2] = 0.2234528123168096
$data[10,3] = 0.2346663685755123
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 3.721411333333334
$data[10,7] = 11.164234
$data[10,8] = 0.1546085275451463
$data[10,9] = 0.1613881693704714
$data[10,10] = 17.52035493851556
$data[10,11] = 157.68319444664
$data[10,12] = 0.03454771028812385
$data[10,13] = 0.03787237563721824
$data[11,0] = 4.707986666666667
$data[11,1] = 14.12396
$data[11,2] = 0.2234528123168096
$data[11,3] = 0.2346663685755123
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 5.194472666666667
$data[11,7] = 15.583418
$data[11,8] = 0.2158078477305768
$data[11,9] = 0.2252711026618443
$data[11,10] = 24.45550805503111
$data[11,11] = 220.09957249528
$data[11,12] = 0.04822287049543519
$data[11,13] = 0.05286355160665641
$data[12,0] = 4.707986666666667
$data[12,1] = 14.12396
$data[12,2] = 0.2234528123168096
$data[12,3] = 0.2346663685755123
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 5.602208666666667
$data[12,7] = 16.806626
$data[12,8] = 0.2327475130727259
$data[12,9] = 0.2429535786722285
$data[12,10] = 26.37512370655111
$data[12,11] = 237.37611335896
$data[12,12] = 0.05200808635584401
$data[12,13] = 0.0570130340394369
$data[13,0] = 4.707986666666667
$data[13,1] = 14.12396
$data[13,2] = 0.2234528123168096
$data[13,3] = 0.2346663685755123
$data[13,4] = 2
$data[13,5] = 1
$data[13,6] = 3.033406
$data[13,7] = 6.066812
$data[13,8] = 0.1260248849423826
$data[13,9] = 0.08770074889103974
$data[13,10] = 14.28123500258667
$data[13,11] = 85.68741001552
$data[13,12] = 0.02816061496227774
$data[13,13] = 0.02058041626361318
$data[14,0] = 4.707986666666667
$data[14,1] = 14.12396
$data[14,2] = 0.2234528123168096
$data[14,3] = 0.2346663685755123
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 6.518398333333334
$data[14,7] = 19.555195
$data[14,8] = 0.2708112267091685
$data[14,9] = 0.2826864004044161
$data[14,10] = 30.68853244135556
$data[14,11] = 276.1967919722
$data[14,12] = 0.06051353021512879
$data[14,13] = 0.06633699102858755
$data[15,0] = 3.020391
$data[15,1] = 6.040782
$data[15,2] = 0.1433553047260926
$data[15,3] = 0.100366212825321
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 3.721411333333334
$data[15,7] = 11.164234
$data[15,8] = 0.1546085275451463
$data[15,9] = 0.1613881693704714
$data[15,10] = 11.240117298498
$data[15,11] = 67.44070379098801
$data[15,12] = 0.02216395257948692
$data[15,13] = 0.01619791935452567
$data[16,0] = 3.020391
$data[16,1] = 6.040782
$data[16,2] = 0.1433553047260926
$data[16,3] = 0.100366212825321
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 5.194472666666667
$data[16,7] = 15.583418
$data[16,8] = 0.2158078477305768
$data[16,9] = 0.2252711026618443
$data[16,10] = 15.689338492146
$data[16,11] = 94.136030952876
$data[16,12] = 0.03093719977369902
$data[16,13] = 0.02260960743315339
$data[17,0] = 3.020391
$data[17,1] = 6.040782
$data[17,2] = 0.1433553047260926
$data[17,3] = 0.100366212825321
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 5.602208666666667
$data[17,7] = 16.806626
$data[17,8] = 0.2327475130727259
$data[17,9] = 0.2429535786722285
$data[17,10] = 16.920860636922
$data[17,11] = 101.525163821532
$data[17,12] = 0.03336559066078084
$data[17,13] = 0.02438433058369024
$data[18,0] = 3.020391
$data[18,1] = 6.040782
$data[18,2] = 0.1433553047260926
$data[18,3] = 0.100366212825321
$data[18,4] = 2
$data[18,5] = 1
$data[18,6] = 3.033406
$data[18,7] = 6.066812
$data[18,8] = 0.1260248849423826
$data[18,9] = 0.08770074889103974
$data[18,10] = 9.162072181746
$data[18,11] = 36.648288726984
$data[18,12] = 0.01806633578398601
$data[18,13] = 0.008802192028138126
$data[19,0] = 3.020391
$data[19,1] = 6.040782
$data[19,2] = 0.1433553047260926
$data[19,3] = 0.100366212825321
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 6.518398333333334
$data[19,7] = 19.555195
$data[19,8] = 0.2708112267091685
$data[19,9] = 0.2826864004044161
$data[19,10] = 19.688111660415
$data[19,11] = 118.12866996249
$data[19,12] = 0.03882222592813978
$data[19,13] = 0.02837216342581352
$data[20,0] = 6.852658666666668
$data[20,1] = 20.557976
$data[20,2] = 0.3252443049075101
$data[20,3] = 0.3415660744707955
$data[20,4] = 3
$data[20,5] = 1
$data[20,6] = 3.721411333333334
$data[20,7] = 11.164234
$data[20,8] = 0.1546085275451463
$data[20,9] = 0.1613881693704714
$data[20,10] = 25.50156162559823
$data[20,11] = 229.514054630384
$data[20,12] = 0.05028554307419472
$data[20,13] = 0.05512472347789978
$data[21,0] = 6.852658666666668
$data[21,1] = 20.557976
$data[21,2] = 0.3252443049075101
$data[21,3] = 0.3415660744707955
$data[21,4] = 3
$data[21,5] = 1
$data[21,6] = 5.194472666666667
$data[21,7] = 15.583418
$data[21,8] = 0.2158078477305768
$data[21,9] = 0.2252711026618443
$data[21,10] = 35.59594813799645
$data[21,11] = 320.3635332419681
$data[21,12] = 0.07019027342871723
$data[21,13] = 0.0769449662279137
$data[22,0] = 6.852658666666668
$data[22,1] = 20.557976
$data[22,2] = 0.3252443049075101
$data[22,3] = 0.3415660744707955
$data[22,4] = 3
$data[22,5] = 1
$data[22,6] = 5.602208666666667
$data[22,7] = 16.806626
$data[22,8] = 0.2327475130727259
$data[22,9] = 0.2429535786722285
$data[22,10] = 38.39002377210845
$data[22,11] = 345.5102139489761
$data[22,12] = 0.07569980310829036
$data[22,13] = 0.08298470014570467
$data[23,0] = 6.852658666666668
$data[23,1] = 20.557976
$data[23,2] = 0.3252443049075101
$data[23,3] = 0.3415660744707955
$data[23,4] = 2
$data[23,5] = 1
$data[23,6] = 3.033406
$data[23,7] = 6.066812
$data[23,8] = 0.1260248849423826
$data[23,9] = 0.08770074889103974
$data[23,10] = 20.78689591541867
$data[23,11] = 124.721375492512
$data[23,12] = 0.04098887610413416
$data[23,13] = 0.02995560052686141
$data[24,0] = 6.852658666666668
$data[24,1] = 20.557976
$data[24,2] = 0.3252443049075101
$data[24,3] = 0.3415660744707955
$data[24,4] = 3
$data[24,5] = 1
$data[24,6] = 6.518398333333334
$data[24,7] = 19.555195
$data[24,8] = 0.2708112267091685
$data[24,9] = 0.2826864004044161
$data[24,10] = 44.66835883170223
$data[24,11] = 402.0152294853201
$data[24,12] = 0.08807980919217363
$data[24,13] = 0.09655608409241588

$ws.Range("G2:T26").Value = $data
Write-Host "Done updating range G2:T26"
